# Add a new example with one PV-bus and one PQ-bus entry (two new time
# entries) to the working-hours log, right before the blank separator /
# summary rows, and let the totals recalculate accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with:
#   row 106: last data row
#   row 107: blank separator row
#   row 108: sum [min]
#   row 109: sum [h]
#   row 110: sum [working weeks]
#
# Insert two fresh rows right before the current blank separator row (107)
# so the two new entries become rows 107 and 108, and the separator /
# summary block shifts down to rows 109-112.
$ws.Rows("107:108").Insert()

# New entry 1: 2014-03-31, 19:30 -> 21:00
$ws.Range("A107").Value = 2014
$ws.Range("B107").Value = 3
$ws.Range("C107").Value = 31
$ws.Range("D107").Value = 0.8125
$ws.Range("E107").Value = 0.875

# New entry 2: 2014-04-01, 14:00 -> 15:30
$ws.Range("A108").Value = 2014
$ws.Range("B108").Value = 4
$ws.Range("C108").Value = 1
$ws.Range("D108").Value = 0.58333333333333337
$ws.Range("E108").Value = 0.64583333333333337

# Fill down the "time spent" formulas for the two new rows, matching the
# pattern used by every other data row in the sheet.
$ws.Range("F107:F108").Formula = "=(E107-D107)*24*60"
$ws.Range("G107:G108").Formula = "=F107/60"

# Update the active selection to reflect where the user ended up editing.
$ws.Range("B109").Select() | Out-Null
